$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.133160666666667
$ws.Range("H2").Value = 24.399482
$ws.Range("I2").Value = 0.3870696756706061
$ws.Range("J2").Value = 0.3870696756706061
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.87733
$ws.Range("N2").Value = 38.63199
$ws.Range("O2").Value = 0.1584922499374361
$ws.Range("P2").Value = 0.1584922499374361
$ws.Range("Q2").Value = 104.7333938476867
$ws.Range("R2").Value = 942.60054462918
$ws.Range("S2").Value = 0.06134754377958804
$ws.Range("T2").Value = 0.06134754377958804
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.133160666666667
$ws.Range("H3").Value = 24.399482
$ws.Range("I3").Value = 0.3870696756706061
$ws.Range("J3").Value = 0.3870696756706061
$ws.Range("O3").Value = 0.4359831802722915
$ws.Range("P3").Value = 0.4359831802722916
$ws.Range("Q3").Value = 288.1024034200391
$ws.Range("R3").Value = 2592.921630780352
$ws.Range("S3").Value = 0.1687558681858353
$ws.Range("T3").Value = 0.1687558681858353
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.133160666666667
$ws.Range("H4").Value = 24.399482
$ws.Range("I4").Value = 0.3870696756706061
$ws.Range("J4").Value = 0.3870696756706061
$ws.Range("M4").Value = 30.51453966666667
$ws.Range("N4").Value = 91.54361900000001
$ws.Range("O4").Value = 0.3755683862706898
$ws.Range("P4").Value = 0.3755683862706898
$ws.Range("Q4").Value = 248.1796537783731
$ws.Range("R4").Value = 2233.616884005358
$ws.Range("S4").Value = 0.1453711334659288
$ws.Range("T4").Value = 0.1453711334659288
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.133160666666667
$ws.Range("H5").Value = 24.399482
$ws.Range("I5").Value = 0.3870696756706061
$ws.Range("J5").Value = 0.3870696756706061
$ws.Range("M5").Value = 2.433908666666667
$ws.Range("N5").Value = 7.301726
$ws.Range("O5").Value = 0.0299561835195825
$ws.Range("P5").Value = 0.0299561835195825
$ws.Range("Q5").Value = 19.79537023399245
$ws.Range("R5").Value = 178.158332105932
$ws.Range("S5").Value = 0.01159513023925395
$ws.Range("T5").Value = 0.01159513023925395
$ws.Range("I6").Value = 0.3559882250904906
$ws.Range("J6").Value = 0.3559882250904906
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 12.87733
$ws.Range("N6").Value = 38.63199
$ws.Range("O6").Value = 0.1584922499374361
$ws.Range("P6").Value = 0.1584922499374361
$ws.Range("Q6").Value = 96.32336844509
$ws.Range("R6").Value = 866.91031600581
$ws.Range("S6").Value = 0.05642137474582631
$ws.Range("T6").Value = 0.0564213747458263
$ws.Range("I7").Value = 0.3559882250904906
$ws.Range("J7").Value = 0.3559882250904906
$ws.Range("O7").Value = 0.4359831802722915
$ws.Range("P7").Value = 0.4359831802722916
$ws.Range("S7").Value = 0.1552048785144405
$ws.Range("T7").Value = 0.1552048785144405
$ws.Range("I8").Value = 0.3559882250904906
$ws.Range("J8").Value = 0.3559882250904906
$ws.Range("M8").Value = 30.51453966666667
$ws.Range("N8").Value = 91.54361900000001
$ws.Range("O8").Value = 0.3755683862706898
$ws.Range("P8").Value = 0.3755683862706898
$ws.Range("Q8").Value = 228.2509842680624
$ws.Range("R8").Value = 2054.258858412561
$ws.Range("S8").Value = 0.1336979232286026
$ws.Range("T8").Value = 0.1336979232286026
$ws.Range("I9").Value = 0.3559882250904906
$ws.Range("J9").Value = 0.3559882250904906
$ws.Range("M9").Value = 2.433908666666667
$ws.Range("N9").Value = 7.301726
$ws.Range("O9").Value = 0.0299561835195825
$ws.Range("P9").Value = 0.0299561835195825
$ws.Range("Q9").Value = 18.20581450199933
$ws.Range("R9").Value = 163.852330517994
$ws.Range("S9").Value = 0.01066404860162118
$ws.Range("T9").Value = 0.01066404860162118
$ws.Range("G10").Value = 5.398902333333333
$ws.Range("H10").Value = 16.196707
$ws.Range("I10").Value = 0.2569420992389033
$ws.Range("J10").Value = 0.2569420992389034
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 12.87733
$ws.Range("N10").Value = 38.63199
$ws.Range("O10").Value = 0.1584922499374361
$ws.Range("P10").Value = 0.1584922499374361
$ws.Range("Q10").Value = 69.52344698410333
$ws.Range("R10").Value = 625.7110228569301
$ws.Range("S10").Value = 0.04072333141202179
$ws.Range("T10").Value = 0.04072333141202179
$ws.Range("G11").Value = 5.398902333333333
$ws.Range("H11").Value = 16.196707
$ws.Range("I11").Value = 0.2569420992389033
$ws.Range("J11").Value = 0.2569420992389034
$ws.Range("O11").Value = 0.4359831802722915
$ws.Range("P11").Value = 0.4359831802722916
$ws.Range("Q11").Value = 191.2462819575502
$ws.Range("R11").Value = 1721.216537617952
$ws.Range("S11").Value = 0.1120224335720158
$ws.Range("T11").Value = 0.1120224335720159
$ws.Range("G12").Value = 5.398902333333333
$ws.Range("H12").Value = 16.196707
$ws.Range("I12").Value = 0.2569420992389033
$ws.Range("J12").Value = 0.2569420992389034
$ws.Range("M12").Value = 30.51453966666667
$ws.Range("N12").Value = 91.54361900000001
$ws.Range("O12").Value = 0.3755683862706898
$ws.Range("P12").Value = 0.3755683862706898
$ws.Range("Q12").Value = 164.7450194069592
$ws.Range("R12").Value = 1482.705174662633
$ws.Range("S12").Value = 0.09649932957615835
$ws.Range("T12").Value = 0.09649932957615838
$ws.Range("G13").Value = 5.398902333333333
$ws.Range("H13").Value = 16.196707
$ws.Range("I13").Value = 0.2569420992389033
$ws.Range("J13").Value = 0.2569420992389034
$ws.Range("M13").Value = 2.433908666666667
$ws.Range("N13").Value = 7.301726
$ws.Range("O13").Value = 0.0299561835195825
$ws.Range("P13").Value = 0.0299561835195825
$ws.Range("Q13").Value = 13.14043517958689
$ws.Range("R13").Value = 118.263916616282
$ws.Range("S13").Value = 0.007697004678707366
$ws.Range("T13").Value = 0.007697004678707368
